# Daily attendance processing - 2026-01-24 21:34:45
# Rotate the "Recorded By" (column G) comma-separated list left by one
# position (move the first entry to the end) for every data row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)   # column G = 7
    $val = $cell.Text

    if ($val -ne $null -and $val -ne "") {
        $parts = $val -split ", "
        if ($parts.Count -gt 1) {
            $rotated = @($parts[1..($parts.Count - 1)]) + @($parts[0])
            $newVal = [string]::Join(", ", $rotated)
            $cell.Value = $newVal
        }
    }
}
